$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipos")

# Update the shared string text used in cell C4: {{item.pendiente}} -> {{item.NumerosSeries}}
$ws.Range("C4").Value = "{{item.NumerosSeries}}"

# Replace the defined name "sucursales" (Equipos!$A$4:$H$5) with "equipos" (Equipos!$A$4:$D$5)
$wb.Names.Item("sucursales").Delete()
$wb.Names.Add("equipos", "=Equipos!`$A`$4:`$D`$5")

# Update the selected/active cell on the sheet from D7 to D4
$ws.Range("D4").Select()
